$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.519.55'
$ws.Range("E2").Value = '  +2.75%  '

$ws.Range("D3").Value = '1.470.79'
$ws.Range("E3").Value = '  +3.69%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9469'
$ws.Range("E5").Value = '  -5.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '281.07'
$ws.Range("E6").Value = '  +2.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3715'
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3187'
$ws.Range("E8").Value = '  +3.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.34'
$ws.Range("E9").Value = '  +3.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.058'
$ws.Range("E10").Value = '  +4.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06675'
$ws.Range("E11").Value = '  +1.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.591'
$ws.Range("E13").Value = '  +3.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.16'
$ws.Range("E14").Value = '  +6.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.250'
$ws.Range("E15").Value = '  +1.45%  '

$ws.Range("D16").Value = '1.477.73'
$ws.Range("E16").Value = '  +3.75%  '

$ws.Range("E17").Value = '  +2.82%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9423'
$ws.Range("E18").Value = '  -6.00%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05746'
$ws.Range("E19").Value = '  -1.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.07'
$ws.Range("E20").Value = '  -3.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.684'
$ws.Range("E21").Value = '  +0.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.69'
$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.22'
$ws.Range("E23").Value = '  +2.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.271'
$ws.Range("E24").Value = '  -2.39%  '

$ws.Range("D25").Value = '20.750.39'
$ws.Range("E25").Value = '  +3.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.295'
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '137.88'
$ws.Range("E27").Value = '  -0.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.59'
$ws.Range("E28").Value = '  +4.00%  '

$ws.Range("D29").Value = '1.642.16'
$ws.Range("E29").Value = '  +3.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.72'
$ws.Range("E30").Value = '  +4.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.954'
$ws.Range("E31").Value = '  +4.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.304'
$ws.Range("E32").Value = '  -2.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8357'
$ws.Range("E33").Value = '  -5.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.621'
$ws.Range("E34").Value = '  +28.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07840'
$ws.Range("E35").Value = '  +1.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06032'
$ws.Range("E36").Value = '  +6.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.927'
$ws.Range("E37").Value = '  +3.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.69'
$ws.Range("E38").Value = '  -4.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02068'
$ws.Range("E39").Value = '  +2.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.125'
$ws.Range("E40").Value = '  +2.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9602'
$ws.Range("E41").Value = '  -4.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1899'
$ws.Range("E42").Value = '  -0.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.451'
$ws.Range("E43").Value = '  -11.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5399'
$ws.Range("E44").Value = '  +1.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.592'
$ws.Range("E45").Value = '  +1.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.38'
$ws.Range("E46").Value = '  +1.39%  '

$ws.Range("E47").Value = '  +11.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5322'
$ws.Range("E48").Value = '  +3.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.827'
$ws.Range("E49").Value = '  +0.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06458'

$ws.Range("E51").Value = '  +0.07%  '
